# Atualização de bases das ligas, do dia: 29-03-2024 às 17:05
#
# The underlying rows got re-matched to the correct fixture/odds data.
# For several pairs (and one triple) of adjacent rows, every column from
# B (id) through AC (PL_AhUnder) needs to be swapped between rows, while
# column A (the running index) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    # Columns B..AC -> 2..29
    return $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 29)).Value2
}

function Set-RowData($row, $data) {
    $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 29)).Value2 = $data
}

# --- Simple pairwise swaps (B:AC), column A untouched ---
$pairs = @(
    @(47, 48),
    @(104, 105),
    @(153, 154)
)

foreach ($p in $pairs) {
    $r1 = $p[0]
    $r2 = $p[1]
    $d1 = Get-RowData $r1
    $d2 = Get-RowData $r2
    Set-RowData $r1 $d2
    Set-RowData $r2 $d1
}

# --- Three-way rotation: 148 <- 150, 149 <- 148, 150 <- 149 (original data) ---
$d148 = Get-RowData 148
$d149 = Get-RowData 149
$d150 = Get-RowData 150

Set-RowData 148 $d150
Set-RowData 149 $d148
Set-RowData 150 $d149
